$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 130803092
$ws.Range("B2").Value = 89193
$ws.Range("E2").Value = 510
$ws.Range("F2").Value = 'Doftskinn'
$ws.Range("G2").Value = 'Cystostereum murrayi'
$ws.Range("H2").Value = '(Berk. & M.A.Curtis.) Pouzar'
$ws.Range("Q2").Value = 424832
$ws.Range("R2").Value = 6712186

# Row 3
$ws.Range("A3").Value = 130803039
$ws.Range("B3").Value = 83223
$ws.Range("E3").Value = 6440
$ws.Range("F3").Value = 'Vitgrynig nållav'
$ws.Range("G3").Value = 'Chaenotheca subroscida'
$ws.Range("H3").Value = '(Eitner) Zahlbr.'
$ws.Range("Q3").Value = 424963
$ws.Range("R3").Value = 6712076

# Row 4
$ws.Range("A4").Value = 130803083
$ws.Range("B4").Value = 83089
$ws.Range("E4").Value = 1312
$ws.Range("F4").Value = 'Gammelgransskål'
$ws.Range("G4").Value = 'Pseudographis pinicola'
$ws.Range("H4").Value = '(Nyl.) Rehm'
$ws.Range("Q4").Value = 424802
$ws.Range("R4").Value = 6712148

# Row 5
$ws.Range("A5").Value = 130803040
$ws.Range("B5").Value = 91828
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = 'Granticka'
$ws.Range("G5").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H5").ClearContents()
$ws.Range("H5").Borders.LineStyle = 0
$ws.Range("Q5").Value = 424793
$ws.Range("R5").Value = 6712247

# Row 8
$ws.Range("A8").Value = 130803042
$ws.Range("B8").Value = 91771
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 5447
$ws.Range("F8").Value = 'Vedticka'
$ws.Range("G8").Value = 'Fuscoporia viticola'
$ws.Range("H8").Value = '(Schwein.) Murrill'
$ws.Range("Q8").Value = 424979
$ws.Range("R8").Value = 6712092

# Row 9
$ws.Range("A9").Value = 130803071
$ws.Range("B9").Value = 91181
$ws.Range("E9").Value = 5685
$ws.Range("F9").Value = 'Gullgröppa'
$ws.Range("G9").Value = 'Pseudomerulius aureus'
$ws.Range("H9").Value = '(Fr.) Jülich'
$ws.Range("Q9").Value = 424873
$ws.Range("R9").Value = 6712126

# Row 10
$ws.Range("A10").Value = 130803064
$ws.Range("B10").Value = 91829
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 5442
$ws.Range("F10").Value = 'Tallticka'
$ws.Range("G10").Value = 'Porodaedalea pini'
$ws.Range("H10").Value = '(Brot.) Murrill'
$ws.Range("Q10").Value = 424893
$ws.Range("R10").Value = 6712101

# Row 11
$ws.Range("A11").Value = 130803067
$ws.Range("B11").Value = 78255
$ws.Range("E11").Value = 228579
$ws.Range("F11").Value = 'Liten svartspik'
$ws.Range("G11").Value = 'Chaenothecopsis nana'
$ws.Range("H11").Value = 'Tibell'
$ws.Range("Q11").Value = 424814
$ws.Range("R11").Value = 6712361

# Row 14
$ws.Range("A14").Value = 130803074
$ws.Range("B14").Value = 79243
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = 'Garnlav'
$ws.Range("G14").Value = 'Alectoria sarmentosa'
$ws.Range("H14").Value = '(Ach.) Ach.'
$ws.Range("Q14").Value = 424801
$ws.Range("R14").Value = 6712403

# Row 15
$ws.Range("A15").Value = 130803041
$ws.Range("B15").Value = 91771
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = 'Vedticka'
$ws.Range("G15").Value = 'Fuscoporia viticola'
$ws.Range("H15").Value = '(Schwein.) Murrill'
$ws.Range("Q15").Value = 424881
$ws.Range("R15").Value = 6712113

# Row 17
$ws.Range("A17").Value = 130803036
$ws.Range("B17").Value = 83223
$ws.Range("E17").Value = 6440
$ws.Range("F17").Value = 'Vitgrynig nållav'
$ws.Range("G17").Value = 'Chaenotheca subroscida'
$ws.Range("H17").Value = '(Eitner) Zahlbr.'
$ws.Range("Q17").Value = 424756
$ws.Range("R17").Value = 6712130

# Row 18
$ws.Range("A18").Value = 130803073
$ws.Range("B18").Value = 79243
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = 'Garnlav'
$ws.Range("G18").Value = 'Alectoria sarmentosa'
$ws.Range("H18").Value = '(Ach.) Ach.'
$ws.Range("Q18").Value = 424873
$ws.Range("R18").Value = 6712251

# Row 19
$ws.Range("A19").Value = 130803084
$ws.Range("B19").Value = 83089
$ws.Range("E19").Value = 1312
$ws.Range("F19").Value = 'Gammelgransskål'
$ws.Range("G19").Value = 'Pseudographis pinicola'
$ws.Range("H19").Value = '(Nyl.) Rehm'
$ws.Range("Q19").Value = 424817
$ws.Range("R19").Value = 6712153

# Row 34
$ws.Range("A34").Value = 130803078
$ws.Range("B34").Value = 79243
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = 'Garnlav'
$ws.Range("G34").Value = 'Alectoria sarmentosa'
$ws.Range("H34").Value = '(Ach.) Ach.'
$ws.Range("Q34").Value = 424951
$ws.Range("R34").Value = 6712137

# Row 36
$ws.Range("A36").Value = 130803082
$ws.Range("B36").Value = 83089
$ws.Range("E36").Value = 1312
$ws.Range("F36").Value = 'Gammelgransskål'
$ws.Range("G36").Value = 'Pseudographis pinicola'
$ws.Range("H36").Value = '(Nyl.) Rehm'
$ws.Range("Q36").Value = 424785
$ws.Range("R36").Value = 6712249

# Row 37
$ws.Range("A37").Value = 130803077
$ws.Range("Q37").Value = 424877
$ws.Range("R37").Value = 6712121

# Row 40
$ws.Range("A40").Value = 130803045
$ws.Range("M40").Value = 'färska spår'
$ws.Range("R40").Value = 6712134

# Row 41
$ws.Range("A41").Value = 130803050
$ws.Range("Q41").Value = 424768
$ws.Range("R41").Value = 6712278

# Row 42
$ws.Range("A42").Value = 130803059
$ws.Range("M42").Value = 'äldre spår'
$ws.Range("Q42").Value = 424858
$ws.Range("R42").Value = 6712137

# Row 45
$ws.Range("A45").Value = 130803046
$ws.Range("M45").Value = 'färska spår'
$ws.Range("Q45").Value = 424910
$ws.Range("R45").Value = 6712215

# Row 46
$ws.Range("A46").Value = 130803058
$ws.Range("M46").Value = 'äldre spår'
$ws.Range("Q46").Value = 424867
$ws.Range("R46").Value = 6712141

# Row 47
$ws.Range("A47").Value = 130803052
$ws.Range("Q47").Value = 424773
$ws.Range("R47").Value = 6712133

# Row 48
$ws.Range("A48").Value = 130803049
$ws.Range("Q48").Value = 424771
$ws.Range("R48").Value = 6712443

# Row 52
$ws.Range("A52").Value = 130803047
$ws.Range("Q52").Value = 424836
$ws.Range("R52").Value = 6712286

# Row 53
$ws.Range("A53").Value = 130803044
$ws.Range("Q53").Value = 424933
$ws.Range("R53").Value = 6712156

# Row 54
$ws.Range("A54").Value = 130803056
$ws.Range("Q54").Value = 424852
$ws.Range("R54").Value = 6712133

# Row 55
$ws.Range("A55").Value = 130803055
$ws.Range("Q55").Value = 424839
$ws.Range("R55").Value = 6712128

# Row 61
$ws.Range("A61").Value = 130848917
$ws.Range("B61").Value = 79001
$ws.Range("E61").Value = 228912
$ws.Range("F61").Value = 'Mörk kolflarnlav'
$ws.Range("G61").Value = 'Carbonicola myrmecina'
$ws.Range("H61").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("K61").ClearContents()
$ws.Range("L61").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("Q61").Value = 424590
$ws.Range("R61").Value = 6712294
$ws.Range("AC61").ClearContents()

# Row 62
$ws.Range("A62").Value = 130848907
$ws.Range("B62").Value = 57884
$ws.Range("E62").Value = 100109
$ws.Range("F62").Value = 'Tretåig hackspett'
$ws.Range("G62").Value = 'Picoides tridactylus'
$ws.Range("H62").Value = '(Linnaeus, 1758)'
$ws.Range("K62").Borders.LineStyle = 0
$ws.Range("L62").Borders.LineStyle = 0
$ws.Range("M62").Value = 'äldre spår'
$ws.Range("N62").Borders.LineStyle = 0
$ws.Range("Q62").Value = 424588
$ws.Range("R62").Value = 6712316
$ws.Range("AC62").Value = 'Ringhack'

# Row 69
$ws.Range("A69").Value = 130848926
$ws.Range("B69").Value = 80308
$ws.Range("D69").Value = 'LC'
$ws.Range("E69").Value = 229497
$ws.Range("F69").Value = 'Korallblylav'
$ws.Range("G69").Value = 'Parmeliella triptophylla'
$ws.Range("H69").Value = '(Ach.) Müll.Arg.'
$ws.Range("Q69").Value = 424630
$ws.Range("R69").Value = 6712465

# Row 70
$ws.Range("A70").Value = 130848920
$ws.Range("B70").Value = 81228
$ws.Range("D70").Value = 'NT'
$ws.Range("E70").Value = 1049
$ws.Range("F70").Value = 'Kortskaftad ärgspik'
$ws.Range("G70").Value = 'Microcalicium ahlneri'
$ws.Range("H70").Value = 'Tibell'
$ws.Range("Q70").Value = 424620
$ws.Range("R70").Value = 6712420
